# Edit: simplify steel description (remove RME)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the text in C2: remove "RME/" from the steel description
$old = $ws.Range("C2").Value2
$new = $old -replace "40% S/LFM\+CDN/RME/H:1", "40% S/LFM+CDN/H:1"
$ws.Range("C2").Value2 = $new

# Apply wrap text to C2 and enlarge row 2 height so the multi-line text is visible
$ws.Range("C2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 320

# Update the selection to match the author's saved selection state
$ws.Range("B2:B5").Select()
